# Adds p-values for the FMB and Shanken-corrected t-stats.
#
# Touches five sheets:
#   - "FMB CAPM"           : add a "p-value" column (E)
#   - "FMB CAPM Shanken"   : insert a new "Alpha" row, add "p-value" column
#   - "FMB FF3F"           : add a "p-value" column (E)
#   - "FMB FF3F Shanken"   : insert a new "Alpha" row, add "p-value" column
#   - "Main Summary"       : mirrors the FMB CAPM / FMB FF3F tables above,
#                            just needs the "p-value" header + values

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: FMB CAPM  (simple extra "p-value" column; no rows are inserted)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FMB CAPM")

$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "p-value"

$ws.Range("E2").Value = 0.000119588677794491
$ws.Range("E3").Value = 0.1161976990350513

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet: FMB CAPM Shanken  (a new "Alpha" row is inserted above the
# existing "Beta_MKT" row, plus the extra "p-value" column)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FMB CAPM Shanken")

$ws.Range("A2").EntireRow.Insert()
$ws.Range("B2:D2").ClearFormats()

$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "Alpha"

$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "p-value"

$ws.Range("B2").Value = 0.01560956353594066
$ws.Range("C2").Value = 0.004066792189690346
$ws.Range("D2").Value = 3.838298788787928
$ws.Range("E2").Value = 0.0001456096855874112

$ws.Range("E3").Value = 0.1209511400843071

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet: FMB FF3F  (simple extra "p-value" column; no rows are inserted)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FMB FF3F")

$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "p-value"

$ws.Range("E2").Value = 0.0009476805940138444
$ws.Range("E3").Value = 0.1249484258211542
$ws.Range("E4").Value = 0.2274980207316721
$ws.Range("E5").Value = 0.00009036051523958477

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet: FMB FF3F Shanken  (a new "Alpha" row is inserted above the
# existing "Beta_MKT" row, plus the extra "p-value" column)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FMB FF3F Shanken")

$ws.Range("A2").EntireRow.Insert()
$ws.Range("B2:D2").ClearFormats()

$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "Alpha"

$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "p-value"

$ws.Range("B2").Value = 0.01129936010281002
$ws.Range("C2").Value = 0.003498952205585338
$ws.Range("D2").Value = 3.229355372380616
$ws.Range("E2").Value = 0.001351352061958799

$ws.Range("E3").Value = 0.1370045312133752
$ws.Range("E4").Value = 0.2421664339711995
$ws.Range("E5").Value = 0.0001467905576828521

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet: Main Summary  (mirrors the FMB CAPM and FMB FF3F tables, just
# adds the "p-value" header + values in column E; no row shifting here)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Main Summary")

$ws.Range("E58").Value = "p-value"

$ws.Range("E60").Value = 0.000119588677794491
$ws.Range("E61").Value = 0.1161976990350513

$ws.Range("E66").Value = "p-value"

$ws.Range("E68").Value = 0.0009476805940138444
$ws.Range("E69").Value = 0.1249484258211542
$ws.Range("E70").Value = 0.2274980207316721
$ws.Range("E71").Value = 0.00009036051523958477

Write-Host "p-value columns added"
